$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 130; this shifts the existing rows 130..217 down
# to 131..218, which matches the diff's downward cascade of D/J/K/L/M/N/P/Q
# values (each old row's data re-appears one row further down, and the
# former last row, 217, becomes the new row 218).
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new record's data. The
# "constant" columns (A,B,C,E,F,G,H,I,N,O,R) repeat the same values used
# throughout this block of rows.
$ws.Range("A130").Value = 8
$ws.Range("B130").Value = "Terminal La Palmera de La Serena"
$ws.Range("C130").Value = "Coquimbo"
$ws.Range("D130").Value = 44762
$ws.Range("E130").Value = 4
$ws.Range("F130").Value = 100112037
$ws.Range("G130").Value = "Cebollín"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 1000
$ws.Range("K130").Value = 1400
$ws.Range("L130").Value = 1600
$ws.Range("M130").Value = 1500
$ws.Range("N130").Value = "$/paquete 6 unidades"
$ws.Range("O130").Value = "Provincia del Elquí"
$ws.Range("P130").Value = 250
$ws.Range("Q130").Value = 6
$ws.Range("R130").Value = "Hortaliza"
